# Reorder the "Recorded By" (column G) names in each attendance row:
# the last two comma-separated entries are swapped, unless the
# second-to-last entry is already the literal "System" (already in the
# desired order, so nothing to do).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Text

    if ([string]::IsNullOrEmpty($val)) {
        continue
    }

    $parts = $val -split ", "
    $n = $parts.Count

    if ($n -ge 2 -and $parts[$n - 2] -ne "System") {
        $tmp = $parts[$n - 1]
        $parts[$n - 1] = $parts[$n - 2]
        $parts[$n - 2] = $tmp
        $cell.Value = ($parts -join ", ")
    }
}
